$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at G (shifts existing G "Labels" column to H)
$ws.Columns("G:G").Insert()

# Header for the newly inserted column G
$ws.Range("G1").Value = "maximaler Histogrammwert"

# "maximaler Histogrammwert" values for rows 2-100 (in row order)
$histValues = @(2024, 9953, 2449, 8822, 2264, 1875, 1663, 1667, 2512, 2513, 4969, 4667, 5133, 5561, 2372, 2657, 1741, 2943, 3181, 5618, 7165, 6791, 7682, 7880, 3401, 3324, 2283, 2756, 4375, 1945, 4170, 3539, 2024, 1580, 3077, 2976, 2281, 3354, 2684, 2321, 2749, 2147, 1510, 1246, 2608, 1652, 7748, 1686, 1959, 2158, 1809, 2255, 5691, 841, 1622, 2343, 2588, 1850, 1848, 1832, 4747, 4135, 2246, 2463, 6457, 2268, 3411, 3820, 4013, 4168, 4165, 5732, 12524, 11407, 3542, 2506, 11423, 2055, 3541, 2600, 2288, 1468, 3334, 924, 999, 1518, 955, 1674, 1606, 2171, 1213, 2231, 2682, 3251, 3801, 2885, 2363, 2680, 5203)

$row = 2
foreach ($v in $histValues) {
    $ws.Cells.Item($row, 7).Value = $v
    $row = $row + 1
}

